$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Duplicate row 39 (the still-blank "TBD" template row) down into rows 40
#    and 41, carrying over both its values and its column formatting
#    (thick/thin left borders used to visually group the columns).
$ws.Range("A39:AY39").Copy()
$ws.Range("A40:AY40").PasteSpecial(-4104)
$ws.Range("A39:AY39").Copy()
$ws.Range("A41:AY41").PasteSpecial(-4104)
$excel.CutCopyMode = $false

# 2) Fill in row 39 ("230910-2") with its real results now that the run
#    finished -- nvidia-smi logging still isn't writing.
$ws.Range("E39").Value = "aborted, nvidia-smi logs are not writing."
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = "NA"
$ws.Range("I39").Value = 1
$ws.Range("J39").Value = "nvidia-smi logs are not writing."
$ws.Range("K39").Value = "NA"
$ws.Range("L39").Value = "NA"
$ws.Range("M39").Value = "NA"
$ws.Range("N39").Value = "NA"
$ws.Range("O39").Value = 0
$ws.Range("R39").Value = 6
$ws.Range("AB39").Value = "NA"
$ws.Range("AW39").Value = 0
$ws.Range("AX39").Value = "NA"
$ws.Range("AY39").Value = "NA"

# 3) Start the new "230910-3" session row (40), then stake out the next
#    "230910-4" template row (41) -- matches the order the two session
#    labels were typed in.
$ws.Range("A40").Value = "230910-3"
$ws.Range("A41").Value = "230910-4"

# 4) Now go back and record row 40's results -- this run had enough memory,
#    but nvidia-smi logging is still broken.
$ws.Range("E40").Value = "aborted, nvidia-smi logs are not writing. But there is enough memory, now."
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = "NA"
$ws.Range("I40").Value = 1
$ws.Range("J40").Value = "nvidia-smi logs are not writing."
$ws.Range("K40").Value = "NA"
$ws.Range("L40").Value = "NA"
$ws.Range("M40").Value = "NA"
$ws.Range("N40").Value = "NA"
$ws.Range("O40").Value = 0
$ws.Range("R40").Value = 6
$ws.Range("AB40").Value = "NA"
$ws.Range("AW40").Value = 0
$ws.Range("AX40").Value = "NA"
$ws.Range("AY40").Value = "NA"

# Restore the formulas in rows 40 & 41 (PasteSpecial of values+formats does
# not carry formulas over, it only carries the cached numbers).
$ws.Range("U40").Formula = "= S40 + T40"
$ws.Range("AE40").Formula = "= 1508.06553301511 + 0.00210606006752809 * (AM40*AN40*AO40) / 5 * U40"
$ws.Range("AQ40").Formula = "= _xlfn.FLOOR.MATH((AJ40 - AM40) / 2)"
$ws.Range("AR40").Formula = "= _xlfn.FLOOR.MATH((AK40 - AN40) / 2)"
$ws.Range("AS40").Formula = "= _xlfn.FLOOR.MATH((AL40 - AO40) / 2)"

$ws.Range("U41").Formula = "= S41 + T41"
$ws.Range("AE41").Formula = "= 1508.06553301511 + 0.00210606006752809 * (AM41*AN41*AO41) / 5 * U41"
$ws.Range("AQ41").Formula = "= _xlfn.FLOOR.MATH((AJ41 - AM41) / 2)"
$ws.Range("AR41").Formula = "= _xlfn.FLOOR.MATH((AK41 - AN41) / 2)"
$ws.Range("AS41").Formula = "= _xlfn.FLOOR.MATH((AL41 - AO41) / 2)"

# 5) Match the on-screen scroll/selection state left behind in the saved
#    file (scrolled up a bit, with H41 selected).
$ws.Activate()
$ws.Range("H41").Select()

Write-Output "done"
